$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Revision #249a7b6: Integer min for rule R30 (row 10) changes from 18 to 1.
$ws.Range("C10").Value = 1
